$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp text (row 1 title, A1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 13:52"

# --- Reorder the small-island-nation rows 180-182 so the country labels read
#     Granada, San Cristobal y Nieves, Seychelles, Groenlandia, Surinam
#     (A180 now shows San Cristobal y Nieves, A181 shows Seychelles, A182 shows Groenlandia) ---
$ws.Range("A180").Value = "San Cristobal y Nieves"
$ws.Range("A181").Value = "Seychelles"
$ws.Range("A182").Value = "Groenlandia"

# --- Updated country statistics ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 367659
$ws.Range("C4").Value = 655
$ws.Range("D4").Value = 19814
$ws.Range("E4").Value = 336902

# Row 21 - Israel
$ws.Range("E21").Value = 8263
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 60

# Row 24 - Australia
$ws.Range("E24").Value = 3313
$ws.Range("G24").Value = 3
$ws.Range("H24").Value = 48

# Row 50 - Grecia
$ws.Range("E50").Value = 1405
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 81

# Row 71 - Bosnia y Herzegovina
$ws.Range("B71").Value = 754
$ws.Range("C71").Value = 80
$ws.Range("E71").Value = 653

# Row 104 - Vietnam
$ws.Range("B104").Value = 249
$ws.Range("C104").Value = 4
$ws.Range("E104").Value = 126

# Row 106 - Montenegro
$ws.Range("D106").Value = 3
$ws.Range("E106").Value = 234
$ws.Range("F106").Value = 7

# Row 180 - San Cristobal y Nieves (after relabel)
$ws.Range("C180").Value = 1

# Row 181 - Seychelles (after relabel)
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 11

# Row 182 - Groenlandia (after relabel)
$ws.Range("B182").Value = 11
$ws.Range("D182").Value = 4
$ws.Range("E182").Value = 7
